# Applies the StructureDefinition metadata refresh (5.0.0 -> 6.0.0, new date,
# Publisher/Jurisdiction replacing the duplicated "Contact" rows) plus the
# Elements-sheet root Short/Definition text update.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" ------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# New IG-publisher run date/time
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 9 used to be "Publisher" with an empty value -> now populated.
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicated "Contact" / "No display for ContactDetail"
# row; retarget it to the new "Jurisdiction" property.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old row 11 ("Contact" / "No display for ContactDetail" duplicate) is
# removed entirely, shifting the remaining rows (old "Description" at row 12,
# etc.) up by one.
$ws.Rows.Item(11).Delete()

# ---- Sheet "Elements" -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# Root extension row: generic placeholder text replaced with the
# profile-specific short/definition text.
$ws2.Range("K2").Value = "Hearing Coverage Indicator"
$ws2.Range("L2").Value = "Indicates whether the member has hearing benefit coverage: Y or N"
